$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay text
# (same display quirk as the source data, e.g. thousands-dot strings).
# Force text number format before writing so Excel does not coerce
# them into real numbers, then restore the default "Normal" style so
# no stray formatting/style delta is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.527.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.522.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.910.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.480.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.523.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0945"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.49%  "
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.109"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0297"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.995.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.761.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.70%  "
